# "Adding the changes we made on may 9th"
#
# accelerometer_selected.xlsx / Sheet1 gains 10 new samples:
#   - two new rows inserted right after the header (new rows 2-3);
#     every previously-existing sample row shifts down by two
#   - eight new rows appended after the old last row (new rows 24-31)
#
# NOTE: this COM shim's Range.Value getter does not return real cell
# contents (it echoes a reflection stub), so the shifted rows cannot be
# "read old, write new cell" - every cell below is written with its final
# literal value directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new rows inserted at the top, right after the header row ---
$ws.Range("A2").Value = 3.577567869787046
$ws.Range("B2").Value = -6.830358225996321
$ws.Range("C2").Value = 7.436059972857905
$ws.Range("A3").Value = 1.296618202114653
$ws.Range("B3").Value = -9.790127622488438
$ws.Range("C3").Value = 8.772651798817325

# --- previously-existing rows, each shifted down by two rows ---
$ws.Range("A4").Value = 2.168524476046066
$ws.Range("B4").Value = -9.669631647141598
$ws.Range("C4").Value = 6.291729929697442
$ws.Range("A5").Value = 2.060296274680465
$ws.Range("B5").Value = -15.79789317115233
$ws.Range("C5").Value = -0.3408509170152625
$ws.Range("A6").Value = -0.5631605959728705
$ws.Range("B6").Value = -8.337644126533947
$ws.Range("C6").Value = 4.529721613088265
$ws.Range("A7").Value = -2.606977947509096
$ws.Range("B7").Value = 2.749581334340653
$ws.Range("C7").Value = 8.704148023826628
$ws.Range("A8").Value = 2.878962026775263
$ws.Range("B8").Value = 1.780854188276234
$ws.Range("C8").Value = 10.20314644976877
$ws.Range("A9").Value = 10.25311217123632
$ws.Range("B9").Value = 0.2482639649954936
$ws.Range("C9").Value = 14.13474797939069
$ws.Range("A10").Value = 17.30286970191231
$ws.Range("B10").Value = -7.345501588852901
$ws.Range("C10").Value = 6.160751110941108
$ws.Range("A11").Value = 16.77036626826363
$ws.Range("B11").Value = -17.12453539200271
$ws.Range("C11").Value = -5.154467277105418
$ws.Range("A12").Value = -10.87548799672837
$ws.Range("B12").Value = -11.36323502050577
$ws.Range("C12").Value = 5.413032257754492
$ws.Range("A13").Value = -2.876455623141776
$ws.Range("B13").Value = -32.54226410586523
$ws.Range("C13").Value = 19.61913488187845
$ws.Range("A14").Value = 50.47583389150513
$ws.Range("B14").Value = -7.503464556530675
$ws.Range("C14").Value = 25.77431648212245
$ws.Range("A15").Value = -1.833547540791238
$ws.Range("B15").Value = -13.35906494519977
$ws.Range("C15").Value = 15.55757014659218
$ws.Range("A16").Value = -26.51655229810853
$ws.Range("B16").Value = 11.14541947380614
$ws.Range("C16").Value = -8.497799986633812
$ws.Range("A17").Value = -6.805932935430175
$ws.Range("B17").Value = 7.032340392223448
$ws.Range("C17").Value = 0.9721723077047082
$ws.Range("A18").Value = -11.96539219176566
$ws.Range("B18").Value = -7.304819296736614
$ws.Range("C18").Value = -2.098955452112834
$ws.Range("A19").Value = 8.602804305145005
$ws.Range("B19").Value = -42.36363804537932
$ws.Range("C19").Value = 16.64096578313493
$ws.Range("A20").Value = 0.5536478593204315
$ws.Range("B20").Value = -8.82534081896356
$ws.Range("C20").Value = 13.45371311672483
$ws.Range("A21").Value = 10.92389786441021
$ws.Range("B21").Value = -9.957679432400008
$ws.Range("C21").Value = 23.96456188391583
$ws.Range("A22").Value = -56.21321548546163
$ws.Range("B22").Value = -8.083060754596882
$ws.Range("C22").Value = -10.20085607707802
$ws.Range("A23").Value = -5.554620584730344
$ws.Range("B23").Value = 4.939894823738217
$ws.Range("C23").Value = -16.05567953204582

# --- new rows appended at the bottom ---
$ws.Range("A24").Value = -5.125987016035118
$ws.Range("B24").Value = -3.317670953866559
$ws.Range("C24").Value = -10.087697013307
$ws.Range("A25").Value = 12.1697812054039
$ws.Range("B25").Value = -11.89667802884434
$ws.Range("C25").Value = -0.7975602624165994
$ws.Range("A26").Value = 13.27590551955934
$ws.Range("B26").Value = 1.014137889798754
$ws.Range("C26").Value = 9.396237328566208
$ws.Range("A27").Value = -4.693403524588408
$ws.Range("B27").Value = -26.48646446354427
$ws.Range("C27").Value = 34.13960077485957
$ws.Range("A28").Value = -13.05436339826206
$ws.Range("B28").Value = -24.43262726167378
$ws.Range("C28").Value = 17.33472581726497
$ws.Range("A29").Value = -15.37497096825698
$ws.Range("B29").Value = 11.57301431993091
$ws.Range("C29").Value = -14.26161232310758
$ws.Range("A30").Value = -9.968023679533061
$ws.Range("B30").Value = -3.754608689092123
$ws.Range("C30").Value = 4.53524044205469
$ws.Range("A31").Value = -3.008851450780502
$ws.Range("B31").Value = -11.21099381420506
$ws.Range("C31").Value = 19.69571330665882

